$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2, shifting the existing data rows down.
$ws.Rows.Item(2).Insert()

# Reset formatting on the newly inserted row back to the default "Normal"
# style (the Insert operation copies the header row's formatting by
# default, which is not what we want for a plain data row).
$ws.Range("A2:T2").Style = "Normal"

# Re-apply the date/time number format used by the other rows' "Fecha"
# column to the new row's D cell.
$ws.Range("D2").NumberFormat = $ws.Range("D3").NumberFormat

# Populate the new row with this week's data.
$ws.Range("A2").Value2 = 1
$ws.Range("B2").Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Range("C2").Value2 = "Arica y Parinacota"
$ws.Range("D2").Value2 = 44699
$ws.Range("E2").Value2 = 15
$ws.Range("F2").Value2 = "Fruta"
$ws.Range("G2").Value2 = 100107
$ws.Range("H2").Value2 = "Otros"
$ws.Range("I2").Value2 = 100107001
$ws.Range("J2").Value2 = "Caqui"
$ws.Range("K2").Value2 = "Mankaki"
$ws.Range("L2").Value2 = "Primera"
$ws.Range("M2").Value2 = 250
$ws.Range("N2").Value2 = 29000
$ws.Range("O2").Value2 = 30000
$ws.Range("P2").Value2 = 29500
$ws.Range("Q2").Value2 = "$/caja 18 kilos granel"
$ws.Range("R2").Value2 = "Región de O'Higgins"
$ws.Range("S2").Value2 = 1639
$ws.Range("T2").Value2 = 18
